# Hma.Calc.xlsx edit: rename the "Index" column to "i" and switch it from a
# 1-based row counter to a 0-based one, narrowing column A to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HMA")

# 1. Header: rename column A's header text from "Index" to "i".
#    (The "testdata" table's first ListColumn name is bound to this cell,
#    so updating the cell updates the table definition too.)
$ws.Range("A1").Value2 = "i"

# 2. Data: re-number rows 2..503 from a 1-based index (1..502) to a
#    0-based index (0..501) by decrementing every existing value by 1.
for ($r = 2; $r -le 503; $r++) {
    $cell = $ws.Range("A$r")
    $cell.Value2 = $cell.Value2 - 1
}

# 3. Narrow column A now that entries are shorter (single/double digits).
$ws.Columns.Item(1).ColumnWidth = 3.14
